# "current state (new derivation)"
# The `gens` sheet's per-generator multiplier column (M) holds a rolling
# window of "current" derivation factors (M35:M67 == 1, i.e. not yet
# derived) followed by the already-derived values (M68:M100). This
# derivation run advances the window: the previously-derived values shift
# up into M35:M67, and M68:M101 reset to the neutral factor of 1 (row 101
# is newly included in the window).

$wb = $excel.ActiveWorkbook

$busses = $wb.Worksheets.Item("busses")
$gens   = $wb.Worksheets.Item("gens")
$lines  = $wb.Worksheets.Item("lines")

# ---------------------------------------------------------------------
# gens!M column: shift the derived factors from M68:M100 up into
# M35:M67, then reset M68:M101 back to 1.
# ---------------------------------------------------------------------
$derived = @(
    1.0210999999999999,
    1.1254,
    1.1141000000000001,
    1.1016999999999999,
    1.0333000000000001,
    1.1564000000000001,
    1.1174999999999999,
    1.085,
    1.1011,
    1.0869,
    1.0125,
    1.0857000000000001,
    1.0584,
    1.0509999999999999,
    1.0422,
    1.0732999999999999,
    0.99990000000000001,
    1.1482000000000001,
    1.1259999999999999,
    1.099,
    1.0983000000000001,
    1.0601,
    1.1509,
    1.0932999999999999,
    1.0643,
    1.1676,
    1.0888,
    1.0705,
    1.0852999999999999,
    1.0576000000000001,
    1.044,
    1.2262999999999999,
    1.1828000000000001
)

for ($i = 0; $i -lt $derived.Length; $i++) {
    $gens.Cells.Item(35 + $i, 13).Value = $derived[$i]
}

for ($row = 68; $row -le 101; $row++) {
    $gens.Cells.Item($row, 13).Value = 1
}

# ---------------------------------------------------------------------
# View / selection state for each sheet, matching the new workbook
# snapshot.
# ---------------------------------------------------------------------

# busses: scrolled up a bit, whole row 50 selected.
[void]$busses.Activate()
[void]$busses.Rows.Item(50).Select()

# lines: no longer the active tab, whole row 13 selected.
[void]$lines.Activate()
[void]$lines.Rows.Item(13).Select()

# gens: becomes the active tab, M68:M101 selected.
[void]$gens.Activate()
[void]$gens.Range("M68:M101").Select()
